$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows: set Status to "Fixed" and Reviser to "Sandro"
$ws.Range("B2").Value = "Fixed"
$ws.Range("C2").Value = "Sandro"

$ws.Range("B5").Value = "Fixed"
$ws.Range("C5").Value = "Sandro"

$ws.Range("B6").Value = "Fixed"
$ws.Range("C6").Value = "Sandro"

$ws.Range("B7").Value = "Fixed"
$ws.Range("C7").Value = "Sandro"

# Add new row 8 for the newly fixed bug
$ws.Range("A8").Value = "Sliding in die animation when killed in air or with knockback"
$ws.Range("B8").Value = "Fixed"
$ws.Range("C8").Value = "Sandro"

# Adjust column A width to fit new content
# (54.15 maps to a stored OOXML width of exactly 55 once Excel rounds
# the character-width figure for the default font)
$ws.Columns.Item(1).ColumnWidth = 54.15

# Update the selected cell in the sheet view
$ws.Range("C5").Select()
